# Edit script: adds a parenthetical red "(This is a change ... )" note after the
# first paragraph's existing text, and appends an extra empty paragraph at the
# very end of the document (just before the final section properties).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Extend the first paragraph ("This is a Microsoft word document.") with:
#       "  "  (two plain spaces, no special formatting)
#       "(This is a change \u2013 Ve"   (red)
#       "rsion for branch alternate"    (red)
#       ")"                              (red)
#    each one becomes its own run, matching the authored edit.
# ---------------------------------------------------------------------------

$firstPara = $d.Paragraphs(1)
$r = $firstPara.Range
# Exclude the paragraph mark from the range so new text lands inside paragraph 1.
$r.End = $r.End - 1
$r.Collapse(0)

# --- run 2: two plain spaces (kept as a distinct run from run 1) -----------
$spacesStart = $r.Start
$r.InsertAfter("  ")
$spacesRange = $d.Range($spacesStart, $spacesStart + 2)
# Re-assigning FormattedText onto itself keeps this text as its own run
# instead of letting the engine silently re-merge it with the preceding,
# identically-formatted run.
$spacesRange.FormattedText = $spacesRange.FormattedText
$r.Collapse(0)

$redColor = 192   # RGB(0xC0,0x00,0x00) -> w:val="C00000"

# --- run 3: "(This is a change \u2013 Ve" ------------------------------------------
$t3 = "(This is a change " + [char]0x2013 + " Ve"
$start3 = $r.End
$r.InsertAfter($t3)
$range3 = $d.Range($start3, $start3 + $t3.Length)
$range3.Font.Color = $redColor
$r.Collapse(0)

# --- run 4: "rsion for branch alternate" ------------------------------------
$t4 = "rsion for branch alternate"
$start4 = $r.End
$r.InsertAfter($t4)
$range4 = $d.Range($start4, $start4 + $t4.Length)
$range4.Font.Color = $redColor
$r.Collapse(0)

# --- run 5: ")" --------------------------------------------------------------
$t5 = ")"
$start5 = $r.End
$r.InsertAfter($t5)
$range5 = $d.Range($start5, $start5 + $t5.Length)
$range5.Font.Color = $redColor
$r.Collapse(0)

# ---------------------------------------------------------------------------
# 2) Append a new, empty paragraph at the very end of the document (right
#    before the sectPr).
# ---------------------------------------------------------------------------

$newPara = $d.Paragraphs.Add()
$newPara.Style = $d.Styles("Normal")

Write-Host "Edit complete"
